# Restore the "R30" rule threshold (column C, row 10 -> cell C10) from 18 to 20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 20
